$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 0.000000000000014210854715202001
$ws.Range("L2").Value = 0.000000000000014210854715202001

$ws.Range("J3").Value = 18.0
$ws.Range("K3").Value = 0.000000000000014210854715202001
$ws.Range("L3").Value = 13.44304095596299

$ws.Range("I4").Value = 18.0
$ws.Range("J4").Value = 27.0
$ws.Range("K4").Value = 15.08191039316677
$ws.Range("L4").Value = 19.52580542088384

$ws.Range("I5").Value = 18.0
$ws.Range("J5").Value = 27.0
$ws.Range("K5").Value = 13.44304095596299
$ws.Range("L5").Value = 19.52580542088384

$ws.Range("I6").Value = 18.0
$ws.Range("J6").Value = 24.0
$ws.Range("K6").Value = 15.63700351079632
$ws.Range("L6").Value = 19.52580542088384

$ws.Range("I7").Value = 27.0
$ws.Range("J7").Value = 33.0
$ws.Range("K7").Value = 19.52580542088384
$ws.Range("L7").Value = 22.32561713860608

$ws.Range("I8").Value = 33.0
$ws.Range("J8").Value = 42.0
$ws.Range("K8").Value = 25.79809015496251
$ws.Range("L8").Value = 29.41820997946633

$ws.Range("I9").Value = 33.0
$ws.Range("J9").Value = 45.0
$ws.Range("K9").Value = 22.32561713860608
$ws.Range("L9").Value = 29.41820997946633

$ws.Range("I10").Value = 33.0
$ws.Range("J10").Value = 42.0
$ws.Range("K10").Value = 25.04554671439056
$ws.Range("L10").Value = 29.41820997946633

$ws.Range("I11").Value = 45.0
$ws.Range("J11").Value = 57.0
$ws.Range("K11").Value = 29.41820997946633
$ws.Range("L11").Value = 38.72905499747877

$ws.Range("I12").Value = 57.0
$ws.Range("J12").Value = 66.0
$ws.Range("K12").Value = 38.72905499747877
$ws.Range("L12").Value = 43.38136781397451

$ws.Range("I13").Value = 57.0
$ws.Range("J13").Value = 63.0
$ws.Range("K13").Value = 39.11351783178653
$ws.Range("L13").Value = 43.38136781397451

$ws.Range("I14").Value = 66.0
$ws.Range("J14").Value = 72.0
$ws.Range("K14").Value = 43.38136781397451
$ws.Range("L14").Value = 46.98092793045322

$ws.Range("I15").Value = 72.0
$ws.Range("J15").Value = 90.0
$ws.Range("K15").Value = 46.98092793045322
$ws.Range("L15").Value = 55.97276790034562

$ws.Range("I16").Value = 90.0
$ws.Range("J16").Value = 120.0
$ws.Range("K16").Value = 55.97276790034562
$ws.Range("L16").Value = 68.04303058184443

$ws.Range("I17").Value = 120.0
$ws.Range("J17").Value = 147.0
$ws.Range("K17").Value = 68.04303058184443
$ws.Range("L17").Value = 85.45746142124123

$ws.Range("I18").Value = 147.0
$ws.Range("J18").Value = 177.0
$ws.Range("K18").Value = 97.72318716379047
$ws.Range("L18").Value = 115.7908874891582

$ws.Range("I19").Value = 147.0
$ws.Range("J19").Value = 192.0
$ws.Range("K19").Value = 85.45746142124123
$ws.Range("L19").Value = 115.7908874891582

$ws.Range("I20").Value = 192.0
$ws.Range("J20").Value = 216.0
$ws.Range("K20").Value = 115.7908874891582
$ws.Range("L20").Value = 126.172413280159

$ws.Range("I21").Value = 216.0
$ws.Range("J21").Value = 225.0
$ws.Range("K21").Value = 129.7778544756666
$ws.Range("L21").Value = 136.1631983179873

$ws.Range("I22").Value = 216.0
$ws.Range("J22").Value = 234.0
$ws.Range("K22").Value = 126.172413280159
$ws.Range("L22").Value = 136.1631983179873

$ws.Range("I23").Value = 216.0
$ws.Range("J23").Value = 225.0
$ws.Range("K23").Value = 130.5428691997917
$ws.Range("L23").Value = 136.1631983179873

$ws.Range("I24").Value = 216.0
$ws.Range("J24").Value = 222.0
$ws.Range("K24").Value = 132.8332923367175
$ws.Range("L24").Value = 136.1631983179873

$ws.Range("I25").Value = 234.0
$ws.Range("J25").Value = 246.0
$ws.Range("K25").Value = 136.1631983179873
$ws.Range("L25").Value = 142.0509972493321

$ws.Range("I26").Value = 246.0
$ws.Range("J26").Value = 258.0
$ws.Range("K26").Value = 142.0509972493321
$ws.Range("L26").Value = 147.7975100160089

$ws.Range("I27").Value = 246.0
$ws.Range("J27").Value = 255.0
$ws.Range("K27").Value = 143.5701975149649
$ws.Range("L27").Value = 147.7975100160089

$ws.Range("I28").Value = 246.0
$ws.Range("J28").Value = 252.0
$ws.Range("K28").Value = 143.2288103544338
$ws.Range("L28").Value = 147.7975100160089

$ws.Range("I29").Value = 246.0
$ws.Range("J29").Value = 252.0
$ws.Range("K29").Value = 143.1905713745582
$ws.Range("L29").Value = 147.7975100160089

$ws.Range("I30").Value = 258.0
$ws.Range("J30").Value = 270.0
$ws.Range("K30").Value = 147.7975100160089
$ws.Range("L30").Value = 154.9364531885381

$ws.Range("I31").Value = 270.0
$ws.Range("J31").Value = 282.0
$ws.Range("K31").Value = 154.9364531885381
$ws.Range("L31").Value = 160.1605965196842

$ws.Range("I32").Value = 270.0
$ws.Range("J32").Value = 276.0
$ws.Range("K32").Value = 155.5476692350085
$ws.Range("L32").Value = 160.1605965196842

$ws.Range("I33").Value = 282.0
$ws.Range("J33").Value = 300.0
$ws.Range("K33").Value = 160.1605965196842
$ws.Range("L33").Value = 172.3199422514753

$ws.Range("I34").Value = 282.0
$ws.Range("J34").Value = 294.0
$ws.Range("K34").Value = 163.6395068759155
$ws.Range("L34").Value = 172.3199422514753

$ws.Range("I35").Value = 282.0
$ws.Range("J35").Value = 288.0
$ws.Range("K35").Value = 167.9939303676939
$ws.Range("L35").Value = 172.3199422514753

$ws.Range("I36").Value = 348.0
$ws.Range("J36").Value = 375.0
$ws.Range("K36").Value = 206.6224847608567
$ws.Range("L36").Value = 226.444825318454

$ws.Range("I37").Value = 300.0
$ws.Range("J37").Value = 348.0
$ws.Range("K37").Value = 172.3199422514753
$ws.Range("L37").Value = 206.6224847608567

$ws.Range("I38").Value = 300.0
$ws.Range("J38").Value = 342.0
$ws.Range("K38").Value = 185.8638278138747
$ws.Range("L38").Value = 206.6224847608567

$ws.Range("I39").Value = 375.0
$ws.Range("J39").Value = 402.0
$ws.Range("K39").Value = 226.444825318454
$ws.Range("L39").Value = 244.7998551128945
$ws.Range("M39").Value = $False

$ws.Range("I40").Value = 402.0
$ws.Range("J40").Value = 444.0
$ws.Range("K40").Value = 244.7998551128945
$ws.Range("L40").Value = 270.7601128539436
$ws.Range("M40").Value = $False

$ws.Range("I41").Value = 402.0
$ws.Range("J41").Value = 426.0
$ws.Range("K41").Value = 256.2211522002381
$ws.Range("L41").Value = 270.7601128539436

$ws.Range("I42").Value = 402.0
$ws.Range("J42").Value = 414.0
$ws.Range("K42").Value = 261.5344267486111
$ws.Range("L42").Value = 270.7601128539436

$ws.Range("I43").Value = 402.0
$ws.Range("J43").Value = 408.0
$ws.Range("K43").Value = 266.3181084792336
$ws.Range("L43").Value = 270.7601128539436

$ws.Range("I44").Value = 402.0
$ws.Range("J44").Value = 408.0
$ws.Range("K44").Value = 267.7386336360777
$ws.Range("L44").Value = 270.7601128539436

$ws.Range("I45").Value = 402.0
$ws.Range("J45").Value = 426.0
$ws.Range("K45").Value = 257.8184581014112
$ws.Range("L45").Value = 270.7601128539436

$ws.Range("I46").Value = 402.0
$ws.Range("J46").Value = 438.0
$ws.Range("K46").Value = 256.7869898520808
$ws.Range("L46").Value = 270.7601128539436

$ws.Range("I47").Value = 444.0
$ws.Range("J47").Value = 453.0
$ws.Range("K47").Value = 270.7601128539436
$ws.Range("L47").Value = 276.8332256961308
$ws.Range("M47").Value = $False

$ws.Range("I48").Value = 453.0
$ws.Range("J48").Value = 471.0
$ws.Range("K48").Value = 276.8332256961308
$ws.Range("L48").Value = 290.5902879800396
$ws.Range("M48").Value = $False

$ws.Range("I49").Value = 453.0
$ws.Range("J49").Value = 465.0
$ws.Range("K49").Value = 282.1390808289829
$ws.Range("L49").Value = 290.5902879800396

$ws.Range("I50").Value = 453.0
$ws.Range("J50").Value = 462.0
$ws.Range("K50").Value = 284.1284706169079
$ws.Range("L50").Value = 290.5902879800396

$ws.Range("I51").Value = 471.0
$ws.Range("J51").Value = 483.0
$ws.Range("K51").Value = 290.5902879800396
$ws.Range("L51").Value = 299.0420988685125
$ws.Range("M51").Value = $False

$ws.Range("I52").Value = 483.0
$ws.Range("J52").Value = 510.0
$ws.Range("K52").Value = 299.0420988685125
$ws.Range("L52").Value = 315.9740768080833
$ws.Range("M52").Value = $False

$ws.Range("I53").Value = 510.0
$ws.Range("J53").Value = 516.0
$ws.Range("K53").Value = 324.5192611168119
$ws.Range("L53").Value = 329.2830629108927

$ws.Range("I54").Value = 510.0
$ws.Range("J54").Value = 528.0
$ws.Range("K54").Value = 315.9740768080833
$ws.Range("L54").Value = 329.2830629108927
$ws.Range("M54").Value = $False

$ws.Range("I55").Value = 510.0
$ws.Range("J55").Value = 516.0
$ws.Range("K55").Value = 325.2030717683239
$ws.Range("L55").Value = 329.2830629108927

$ws.Range("I56").Value = 528.0
$ws.Range("J56").Value = 534.0
$ws.Range("K56").Value = 329.5299216078971
$ws.Range("L56").Value = 333.4440938846141

$ws.Range("I57").Value = 528.0
$ws.Range("J57").Value = 534.0
$ws.Range("K57").Value = 329.2830629108927
$ws.Range("L57").Value = 333.4440938846141
$ws.Range("M57").Value = $False

$ws.Range("I58").Value = 534.0
$ws.Range("J58").Value = 540.0
$ws.Range("K58").Value = 336.0659134034996
$ws.Range("L58").Value = 340.3465076658096

$ws.Range("I59").Value = 534.0
$ws.Range("J59").Value = 543.0
$ws.Range("K59").Value = 333.4440938846141
$ws.Range("L59").Value = 340.3465076658096
$ws.Range("M59").Value = $False

$ws.Range("I60").Value = 534.0
$ws.Range("J60").Value = 540.0
$ws.Range("K60").Value = 336.8593551736687
$ws.Range("L60").Value = 340.3465076658096

$ws.Range("I61").Value = 543.0
$ws.Range("J61").Value = 549.0
$ws.Range("K61").Value = 340.3465076658096
$ws.Range("L61").Value = 345.0811463029151
$ws.Range("M61").Value = $False

$ws.Range("I62").Value = 549.0
$ws.Range("J62").Value = 555.0
$ws.Range("K62").Value = 345.0811463029151
$ws.Range("L62").Value = 349.2397617328445
$ws.Range("M62").Value = $False

$ws.Range("I63").Value = 555.0
$ws.Range("J63").Value = 556.0
$ws.Range("K63").Value = 349.2397617328445
$ws.Range("L63").Value = 350.2397617328445
$ws.Range("M63").Value = $False

$ws.Range("I64").Value = 556.0
$ws.Range("J64").Value = 556.0
$ws.Range("K64").Value = 350.2397617328445
$ws.Range("L64").Value = 350.2397617328445
$ws.Range("M64").Value = $False

Write-Output "Applied gate/critical-path value updates"
